# Diagnostic sheet: "disconnected_elements" contingency cell block.
#   B1 = 0              (bold, thin boxed border, centered/top-aligned)
#   A2 = 0              (same formatting as B1)
#   B2 = "disconnected_elements"   (plain, shared string)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Build the boxed/bold/centered style once on B1 ...
$r1 = $ws.Range("B1")
$r1.Font.Bold = $true
$r1.HorizontalAlignment = -4108  # xlCenter
$r1.VerticalAlignment = -4160    # xlTop
$r1.Borders.LineStyle = 1        # xlContinuous (thin box around the cell)

# ... then copy/paste just the formatting onto A2 so both cells resolve to
# the very same cell-style record instead of minting a second, equivalent one.
$r2 = $ws.Range("A2")
$r1.Copy()
$r2.PasteSpecial(-4122)          # xlPasteFormats
